# Fixed a bug in respin stats
# The underlying data-generation logic for the per-symbol respin statistics
# table was reordered; this updates the worksheet values to reflect the
# corrected row ordering/values (A2:F25), leaving headers (row 1) and the
# totals row (row 26) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1201, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(701, 3, 90, 45, 97, 15),
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(801, 3, 67, 65, 52, 45),
    @(101, 9, 30, 15, 60, 15),
    @(301, 6, 45, 30, 60, 45),
    @(401, 9, 48, 67, 75, 45),
    @(901, 16, 15, 45, 60, 60),
    @(902, 1, 0, 0, 0, 0),
    @(1001, 18, 30, 75, 60, 72),
    @(501, 9, 52, 30, 75, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(2, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(1101, 0, 15, 30, 30, 0),
    @(1, 0, 2, 2, 2, 2),
    @(502, 0, 4, 0, 0, 0),
    @(802, 0, 4, 5, 4, 0),
    @(602, 0, 0, 4, 0, 9),
    @(402, 0, 0, 4, 0, 0),
    @(702, 0, 0, 0, 4, 0),
    @(1002, 0, 0, 0, 0, 9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
